$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates (odds changed) ---
$ws.Cells.Item(3, 7).Value = 3.1
$ws.Cells.Item(3, 9).Value = 2.25
$ws.Cells.Item(3, 11).Value = 2.1
$ws.Cells.Item(3, 24).Value = 15
$ws.Cells.Item(3, 25).Value = 11
$ws.Cells.Item(3, 27).Value = 23
$ws.Cells.Item(3, 29).Value = 10
$ws.Cells.Item(3, 38).Value = 19
$ws.Cells.Item(3, 50).Value = 13

# --- Row 4 updates (odds changed) ---
$ws.Cells.Item(4, 8).Value = 3.5
$ws.Cells.Item(4, 9).Value = 3.2
$ws.Cells.Item(4, 11).Value = 2.25
$ws.Cells.Item(4, 21).Value = 1.62
$ws.Cells.Item(4, 22).Value = 2.2
$ws.Cells.Item(4, 23).Value = 9.5
$ws.Cells.Item(4, 33).Value = 151
$ws.Cells.Item(4, 54).Value = 151

# --- Insert two new rows at position 6 for the two new SPAIN - LALIGA2 matches ---
$ws.Range("A6:A7").EntireRow.Insert()

# --- New row 6: Cadiz CF vs Cordoba ---
$ws.Cells.Item(6, 1).Value = "4vIbqhgd"
$ws.Cells.Item(6, 2).Value = "17/11/2024"
$ws.Cells.Item(6, 3).Value = "12:15"
$ws.Cells.Item(6, 4).Value = "SPAIN - LALIGA2"
$ws.Cells.Item(6, 5).Value = "Cadiz CF"
$ws.Cells.Item(6, 6).Value = "Cordoba"
$ws.Cells.Item(6, 7).Value = 1.95
$ws.Cells.Item(6, 8).Value = 3.3
$ws.Cells.Item(6, 9).Value = 4
$ws.Cells.Item(6, 10).Value = 2.63
$ws.Cells.Item(6, 11).Value = 2.1
$ws.Cells.Item(6, 12).Value = 4.5
$ws.Cells.Item(6, 13).Value = 1.06
$ws.Cells.Item(6, 14).Value = 10
$ws.Cells.Item(6, 15).Value = 1.33
$ws.Cells.Item(6, 16).Value = 3.25
$ws.Cells.Item(6, 17).Value = 2.05
$ws.Cells.Item(6, 18).Value = 1.75
$ws.Cells.Item(6, 19).Value = 1.44
$ws.Cells.Item(6, 20).Value = 2.63
$ws.Cells.Item(6, 21).Value = 1.83
$ws.Cells.Item(6, 22).Value = 1.83
$ws.Cells.Item(6, 23).Value = 7
$ws.Cells.Item(6, 24).Value = 9
$ws.Cells.Item(6, 25).Value = 9
$ws.Cells.Item(6, 26).Value = 17
$ws.Cells.Item(6, 27).Value = 17
$ws.Cells.Item(6, 28).Value = 29
$ws.Cells.Item(6, 29).Value = 9
$ws.Cells.Item(6, 30).Value = 6.5
$ws.Cells.Item(6, 31).Value = 15
$ws.Cells.Item(6, 32).Value = 51
$ws.Cells.Item(6, 33).Value = 301
$ws.Cells.Item(6, 34).Value = 11
$ws.Cells.Item(6, 35).Value = 19
$ws.Cells.Item(6, 36).Value = 13
$ws.Cells.Item(6, 37).Value = 41
$ws.Cells.Item(6, 38).Value = 34
$ws.Cells.Item(6, 39).Value = 41
$ws.Cells.Item(6, 40).Value = 4
$ws.Cells.Item(6, 41).Value = 11
$ws.Cells.Item(6, 42).Value = 23
$ws.Cells.Item(6, 43).Value = 41
$ws.Cells.Item(6, 44).Value = 51
$ws.Cells.Item(6, 45).Value = 151
$ws.Cells.Item(6, 46).Value = 2.63
$ws.Cells.Item(6, 47).Value = 8.5
$ws.Cells.Item(6, 48).Value = 51
$ws.Cells.Item(6, 49).Value = 5.5
$ws.Cells.Item(6, 50).Value = 21
$ws.Cells.Item(6, 51).Value = 29
$ws.Cells.Item(6, 52).Value = 81
$ws.Cells.Item(6, 53).Value = 101
$ws.Cells.Item(6, 54).Value = 251
$ws.Cells.Item(6, 55).Value = 81
$ws.Cells.Item(6, 56).Value = 81

# --- New row 7: R. Oviedo vs Tenerife ---
$ws.Cells.Item(7, 1).Value = "b1wYci8k"
$ws.Cells.Item(7, 2).Value = "17/11/2024"
$ws.Cells.Item(7, 3).Value = "12:15"
$ws.Cells.Item(7, 4).Value = "SPAIN - LALIGA2"
$ws.Cells.Item(7, 5).Value = "R. Oviedo"
$ws.Cells.Item(7, 6).Value = "Tenerife"
$ws.Cells.Item(7, 7).Value = 1.67
$ws.Cells.Item(7, 8).Value = 3.5
$ws.Cells.Item(7, 9).Value = 5.25
$ws.Cells.Item(7, 10).Value = 2.4
$ws.Cells.Item(7, 11).Value = 1.95
$ws.Cells.Item(7, 12).Value = 6.5
$ws.Cells.Item(7, 13).Value = 1.1
$ws.Cells.Item(7, 14).Value = 7
$ws.Cells.Item(7, 15).Value = 1.5
$ws.Cells.Item(7, 16).Value = 2.5
$ws.Cells.Item(7, 17).Value = 2.5
$ws.Cells.Item(7, 18).Value = 1.5
$ws.Cells.Item(7, 19).Value = 1.57
$ws.Cells.Item(7, 20).Value = 2.25
$ws.Cells.Item(7, 21).Value = 2.5
$ws.Cells.Item(7, 22).Value = 1.5
$ws.Cells.Item(7, 23).Value = 5
$ws.Cells.Item(7, 24).Value = 6.5
$ws.Cells.Item(7, 25).Value = 9.5
$ws.Cells.Item(7, 26).Value = 12
$ws.Cells.Item(7, 27).Value = 17
$ws.Cells.Item(7, 28).Value = 41
$ws.Cells.Item(7, 29).Value = 6.5
$ws.Cells.Item(7, 30).Value = 7
$ws.Cells.Item(7, 31).Value = 23
$ws.Cells.Item(7, 32).Value = 101
$ws.Cells.Item(7, 33).Value = 201
$ws.Cells.Item(7, 34).Value = 10
$ws.Cells.Item(7, 35).Value = 26
$ws.Cells.Item(7, 36).Value = 19
$ws.Cells.Item(7, 37).Value = 67
$ws.Cells.Item(7, 38).Value = 51
$ws.Cells.Item(7, 39).Value = 51
$ws.Cells.Item(7, 40).Value = 3.4
$ws.Cells.Item(7, 41).Value = 9.5
$ws.Cells.Item(7, 42).Value = 26
$ws.Cells.Item(7, 43).Value = 34
$ws.Cells.Item(7, 44).Value = 67
$ws.Cells.Item(7, 45).Value = 301
$ws.Cells.Item(7, 46).Value = 2.25
$ws.Cells.Item(7, 47).Value = 10
$ws.Cells.Item(7, 48).Value = 81
$ws.Cells.Item(7, 49).Value = 7
$ws.Cells.Item(7, 50).Value = 34
$ws.Cells.Item(7, 51).Value = 41
$ws.Cells.Item(7, 52).Value = 151
$ws.Cells.Item(7, 53).Value = 201
$ws.Cells.Item(7, 54).Value = 501
$ws.Cells.Item(7, 55).Value = 81
$ws.Cells.Item(7, 56).Value = 81
